$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.929.75'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.642.98'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.61'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5053'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2580'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06401'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.61'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07778'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.284'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '1.620.42'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5439'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '0.0₅7879'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.94'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").Value = '25.978.62'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '198.27'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.414'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.973'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.002'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.867'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.99'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1143'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.879'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.75'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05017'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.271'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.201'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.375'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8945'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.613'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '1.145.39'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5553'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01563'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.696'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8241'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.03'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  +9.02%  '
$ws.Range("D45").Value = '1.783.40'
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4528'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.44'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05062'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09534'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.48%  '
